# The source data accidentally had an extra leading column (old column A,
# containing the GENE numbers 3/5/7/13) duplicated at the end in column F.
# Remove the stray leading column so the table starts cleanly at column A
# again: everything shifts one column to the left (B->A, C->B, D->C, E->D,
# F->E) and the old column A is discarded, matching the final result layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A:A").Delete()
